{"js": "// Update the division-problem worksheet cells with new values, matching\n// the old text exactly so formatting (font, size, run properties) is kept.\nconst replacements = [\n  [\"19\u00f79=\", \"70\u00f73=\"],\n  [\"20\u00f75=\", \"73\u00f73=\"],\n  [\"62\u00f79=\", \"21\u00f75=\"],\n  [\"97\u00f77=\", \"58\u00f74=\"],\n  [\"24\u00f78=\", \"16\u00f79=\"],\n  [\"70\u00f72=\", \"98\u00f76=\"],\n  [\"82\u00f72=\", \"70\u00f78=\"],\n  [\"16\u00f78=\", \"92\u00f76=\"],\n  [\"61\u00f79=\", \"16\u00f79=\"],\n  [\"85\u00f77=\", \"82\u00f74=\"],\n  [\"79\u00f75=\", \"40\u00f77=\"],\n  [\"16\u00f75=\", \"41\u00f77=\"],\n  [\"30\u00f73=\", \"30\u00f78=\"],\n  [\"78\u00f78=\", \"48\u00f76=\"],\n  [\"97\u00f75=\", \"10\u00f72=\"],\n  [\"80\u00f79=\", \"51\u00f72=\"],\n  [\"93\u00f76=\", \"35\u00f77=\"],\n  [\"51\u00f73=\", \"11\u00f79=\"],\n  [\"91\u00f76=\", \"27\u00f73=\"],\n  [\"42\u00f74=\", \"54\u00f72=\"],\n  [\"98\u00f74=\", \"11\u00f79=\"],\n  [\"34\u00f76=\", \"22\u00f78=\"],\n  [\"73\u00f78=\", \"23\u00f72=\"],\n  [\"63\u00f73=\", \"88\u00f77=\"],\n  [\"81\u00f78=\", \"75\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem worksheet cells with new values, matching\n# the old text exactly so formatting (font, size, run properties) is kept.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"19\u00f79=\", \"70\u00f73=\"),\n  @(\"20\u00f75=\", \"73\u00f73=\"),\n  @(\"62\u00f79=\", \"21\u00f75=\"),\n  @(\"97\u00f77=\", \"58\u00f74=\"),\n  @(\"24\u00f78=\", \"16\u00f79=\"),\n  @(\"70\u00f72=\", \"98\u00f76=\"),\n  @(\"82\u00f72=\", \"70\u00f78=\"),\n  @(\"16\u00f78=\", \"92\u00f76=\"),\n  @(\"61\u00f79=\", \"16\u00f79=\"),\n  @(\"85\u00f77=\", \"82\u00f74=\"),\n  @(\"79\u00f75=\", \"40\u00f77=\"),\n  @(\"16\u00f75=\", \"41\u00f77=\"),\n  @(\"30\u00f73=\", \"30\u00f78=\"),\n  @(\"78\u00f78=\", \"48\u00f76=\"),\n  @(\"97\u00f75=\", \"10\u00f72=\"),\n  @(\"80\u00f79=\", \"51\u00f72=\"),\n  @(\"93\u00f76=\", \"35\u00f77=\"),\n  @(\"51\u00f73=\", \"11\u00f79=\"),\n  @(\"91\u00f76=\", \"27\u00f73=\"),\n  @(\"42\u00f74=\", \"54\u00f72=\"),\n  @(\"98\u00f74=\", \"11\u00f79=\"),\n  @(\"34\u00f76=\", \"22\u00f78=\"),\n  @(\"73\u00f78=\", \"23\u00f72=\"),\n  @(\"63\u00f73=\", \"88\u00f77=\"),\n  @(\"81\u00f78=\", \"75\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
